$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 25: new "Notification for EOHS Input during Combined Check" entry ---
# A25: sequence number (row 24 has 23, so this is 24)
$ws.Range("A25").Value = 24

# B25: title/description
$ws.Range("B25").Value = "Notification for EOHS Input during Combined Check"

# C25: related objects/scripts (two lines)
$ws.Range("C25").Value = "operaDatalog_SX" + [char]10 + "Script_Send_EOHS_Message"

# D25 / E25: same "Using" / "Taki" values as the row above (D24/E24)
$ws.Range("D25").Value = $ws.Range("D24").Value2
$ws.Range("E25").Value = $ws.Range("E24").Value2

# F25: hyperlink formula to the new document
$ws.Range("F25").Formula = '=HYPERLINK(".\Notification%20for%20EOHS%20Input%20Inspection.docx", "Notification for EOHS Input Inspection")'

# Match the hyperlink cell style used by the existing rows (e.g. F24) so the
# new link cell renders the same way (centered, hyperlink font/border).
$ws.Range("F24").Copy()
$ws.Range("F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the active view/selection to reflect where the edit happened ---
$ws.Range("A22").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("F27").Select()
